# Apply profit-table updates to Aegis_Profits sheets (scheduled runner refresh).
# Values below come from an upstream price-data resync; only the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 573.37
$ws.Range("I15").Value = 573.37
$ws.Range("K15").Value = 1720.11
$ws.Range("M15").Value = -1551.11

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 551.37933
$ws.Range("I28").Value = 291.45834
$ws.Range("J28").Value = 1799
$ws.Range("K28").Value = 291.45834
$ws.Range("L28").Value = 1799
$ws.Range("M28").Value = 193.54166
$ws.Range("N28").Value = -2769

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 3333.3333
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6248

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 3333.3333
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -31240

# Row 69: Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 3428
$ws.Range("J69").Value = 3272.5
$ws.Range("L69").Value = 9817.5
$ws.Range("N69").Value = -11565.5

# Row 72: Surgical Substitution (L)
$ws.Range("H72").Value = 3428
$ws.Range("J72").Value = 3272.5
$ws.Range("L72").Value = 29452.5
$ws.Range("N72").Value = -38188.5

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1048.1538
$ws.Range("J112").Value = 1094.6666
$ws.Range("L112").Value = 3283.9998
$ws.Range("N112").Value = -5499.9998

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1172.75
$ws.Range("I137").Value = 1139.5714
$ws.Range("J137").Value = 1272.2858
$ws.Range("K137").Value = 3418.7142
$ws.Range("L137").Value = 3816.8574
$ws.Range("M137").Value = -868.7142000000003
$ws.Range("N137").Value = -8916.857400000001

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1571.4923
$ws.Range("I138").Value = 1441.3
$ws.Range("J138").Value = 1629.3556
$ws.Range("K138").Value = 4323.9
$ws.Range("L138").Value = 4888.066800000001
$ws.Range("M138").Value = 816.1000000000004
$ws.Range("N138").Value = -15168.0668


$ws = $wb.Worksheets.Item("ARM")
# Row 6: Don't Hit Me One More Time
$ws.Range("H6").Value = 20938.25
$ws.Range("I6").Value = 37876.5
$ws.Range("J6").Value = 4000
$ws.Range("K6").Value = 37876.5
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = -37703.5
$ws.Range("N6").Value = -4346

# Row 46: Get Me the Usual
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2118.9
$ws.Range("I61").Value = 1659.2778
$ws.Range("K61").Value = 1659.2778
$ws.Range("M61").Value = -1447.2778

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2118.9
$ws.Range("I136").Value = 1659.2778
$ws.Range("K136").Value = 4977.8334
$ws.Range("M136").Value = -2427.8334


$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 440
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 680
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 680
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -1026

# Row 80: Unbreaker
$ws.Range("H80").Value = 1932.9286
$ws.Range("I80").Value = 897.5714
$ws.Range("K80").Value = 897.5714
$ws.Range("M80").Value = 100.4286

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 1932.9286
$ws.Range("I83").Value = 897.5714
$ws.Range("K83").Value = 4487.857
$ws.Range("M83").Value = 504.143

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2247.7778
$ws.Range("I99").Value = 2293.3333
$ws.Range("J99").Value = 2238.6667
$ws.Range("K99").Value = 2293.3333
$ws.Range("L99").Value = 2238.6667
$ws.Range("M99").Value = -795.3332999999998
$ws.Range("N99").Value = -5234.6667

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 78732.92
$ws.Range("I105").Value = 54506.58
$ws.Range("K105").Value = 54506.58
$ws.Range("M105").Value = -52759.58


$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 4011.4707
$ws.Range("I86").Value = 3700
$ws.Range("K86").Value = 3700
$ws.Range("M86").Value = -2577

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 4011.4707
$ws.Range("I89").Value = 3700
$ws.Range("K89").Value = 18500
$ws.Range("M89").Value = -12884

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1489.2142
$ws.Range("I105").Value = 1295.5714
$ws.Range("J105").Value = 1682.8572
$ws.Range("K105").Value = 1295.5714
$ws.Range("L105").Value = 1682.8572
$ws.Range("M105").Value = 451.4286
$ws.Range("N105").Value = -5176.8572


$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 852.5
$ws.Range("J131").Value = 867.52576
$ws.Range("L131").Value = 2602.57728
$ws.Range("N131").Value = -12682.57728


$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me
$ws.Range("H5").Value = 1669666.6
$ws.Range("I5").Value = 5000000
$ws.Range("K5").Value = 5000000
$ws.Range("M5").Value = -4999888

# Row 74: The Unfortunate Retirony
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77: Life Ends at Retirement (L)
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 1385.1428
$ws.Range("I113").Value = 1218.6
$ws.Range("J113").Value = 1477.6666
$ws.Range("K113").Value = 1218.6
$ws.Range("L113").Value = 1477.6666
$ws.Range("M113").Value = 951.4000000000001
$ws.Range("N113").Value = -5817.6666


$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 229750
$ws.Range("J2").Value = 8636.362999999999
$ws.Range("L2").Value = 8636.362999999999
$ws.Range("N2").Value = -8860.362999999999

# Row 16: Saddle Sore
$ws.Range("H16").Value = 424477.47
$ws.Range("I16").Value = 84332.586
$ws.Range("J16").Value = 716030.2
$ws.Range("K16").Value = 84332.586
$ws.Range("L16").Value = 716030.2
$ws.Range("M16").Value = -84162.586
$ws.Range("N16").Value = -716370.2

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 9388.666999999999
$ws.Range("I132").Value = 16333.333
$ws.Range("J132").Value = 5916.3335
$ws.Range("K132").Value = 48999.999
$ws.Range("L132").Value = 17749.0005
$ws.Range("M132").Value = -46469.999
$ws.Range("N132").Value = -22809.0005

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1504.5555
$ws.Range("I136").Value = 1268.0358
$ws.Range("J136").Value = 2332.375
$ws.Range("K136").Value = 3804.1074
$ws.Range("L136").Value = 6997.125
$ws.Range("M136").Value = -1254.1074
$ws.Range("N136").Value = -12097.125


$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Range("H2").Value = 9333
$ws.Range("I2").Value = 9333
$ws.Range("K2").Value = 9333
$ws.Range("M2").Value = -9221

# Row 75: Storm upon Bald Mountain
$ws.Range("H75").Value = 28276.666
$ws.Range("I75").Value = 15000
$ws.Range("K75").Value = 15000
$ws.Range("M75").Value = -14064

# Row 78: Abrupt Apprentices (L)
$ws.Range("H78").Value = 28276.666
$ws.Range("I78").Value = 15000
$ws.Range("K78").Value = 45000
$ws.Range("M78").Value = -40320

# Row 80: Healing with Flair
$ws.Range("H80").Value = 39990
$ws.Range("J80").Value = 39990
$ws.Range("L80").Value = 39990
$ws.Range("N80").Value = -41986

# Row 83: Pants Fit for Battle (L)
$ws.Range("H83").Value = 39990
$ws.Range("J83").Value = 39990
$ws.Range("L83").Value = 119970
$ws.Range("N83").Value = -129954

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2270.7273
$ws.Range("I132").Value = 2296.7693
$ws.Range("J132").Value = 2207.25
$ws.Range("K132").Value = 6890.3079
$ws.Range("L132").Value = 6621.75
$ws.Range("M132").Value = -4360.3079
$ws.Range("N132").Value = -11681.75
